$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '66.215.81'
$c.Style = "Normal"
$ws.Range('E2').Value = '  -0.15%  '
$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '3.565.81'
$c.Style = "Normal"
$ws.Range('E3').Value = '  +1.14%  '
$ws.Range('E4').Value = '  +0.06%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '606.48'
$c.Style = "Normal"
$ws.Range('E5').Value = '  -0.06%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '144.22'
$c.Style = "Normal"
$ws.Range('E6').Value = '  -0.68%  '
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '3.563.07'
$c.Style = "Normal"
$ws.Range('E7').Value = '  +1.10%  '
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('E9').Value = '  +2.48%  '
$ws.Range('E10').Value = '  -0.17%  '
$ws.Range('E11').Value = '  -2.70%  '
$ws.Range('E12').Value = '  -0.17%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '4.171.53'
$c.Style = "Normal"
$ws.Range('E13').Value = '  +1.21%  '
$ws.Range('E14').Value = '  -0.36%  '
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '30.36'
$c.Style = "Normal"
$ws.Range('E15').Value = '  -0.22%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '3.559.83'
$c.Style = "Normal"
$ws.Range('E16').Value = '  +1.10%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '66.278.51'
$c.Style = "Normal"
$ws.Range('E17').Value = '  -0.13%  '
$ws.Range('E19').Value = '  +6.47%  '
$ws.Range('E20').Value = '  +0.16%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '14.79'
$c.Style = "Normal"
$ws.Range('E21').Value = '  -1.04%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '431.23'
$c.Style = "Normal"
$ws.Range('E22').Value = '  +0.96%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '0.614'
$c.Style = "Normal"
$ws.Range('E23').Value = '  +2.14%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '79.47'
$c.Style = "Normal"
$ws.Range('E24').Value = '  +1.60%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '3.709.41'
$c.Style = "Normal"
$ws.Range('E25').Value = '  +1.36%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('E27').Value = '  -2.08%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '9.17'
$c.Style = "Normal"
$ws.Range('E28').Value = '  -1.37%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '2.50'
$c.Style = "Normal"
$ws.Range('E29').Value = '  +0.78%  '
$ws.Range('E30').Value = '  -1.59%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '0.994'
$c.Style = "Normal"
$ws.Range('E31').Value = '  -0.57%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '3.562.78'
$c.Style = "Normal"
$ws.Range('E32').Value = '  +1.57%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '25.44'
$c.Style = "Normal"
$ws.Range('E33').Value = '  +0.62%  '
$ws.Range('E34').Value = '  -2.01%  '
$ws.Range('E35').Value = '  -8.23%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '7.88'
$c.Style = "Normal"
$ws.Range('E36').Value = '  +0.99%  '
$ws.Range('E37').Value = '  +0.03%  '
$ws.Range('E38').Value = '  -0.94%  '
$ws.Range('E39').Value = '  -0.14%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '173.30'
$c.Style = "Normal"
$ws.Range('E40').Value = '  +1.49%  '
$ws.Range('E41').Value = '  -0.93%  '
$ws.Range('E42').Value = '  +0.33%  '
$ws.Range('E43').Value = '  -0.08%  '
$ws.Range('E44').Value = '  +1.91%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '45.95'
$c.Style = "Normal"
$ws.Range('E45').Value = '  +1.06%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('E47').Value = '  -1.29%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '2.44'
$c.Style = "Normal"
$ws.Range('E48').Value = '  +1.08%  '
$ws.Range('E49').Value = '  -3.96%  '
$ws.Range('E50').Value = '  -0.65%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '23.40'
$c.Style = "Normal"
$ws.Range('E51').Value = '  +4.83%  '
